$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.008.50"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.881.02"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'243.11"
$ws.Range("E5").Value = "  -3.62%  "
$ws.Range("D6").Value = "'0.9984"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.4963"
$ws.Range("E7").Value = "  -3.25%  "
$ws.Range("D8").Value = "'0.2924"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "'0.06649"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "1.880.87"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").Value = "'16.76"
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("D12").Value = "'0.07254"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "'0.6700"
$ws.Range("E13").Value = "  -3.40%  "
$ws.Range("D14").Value = "'86.59"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "'4.880"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "29.991.79"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "'0.000007905"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").Value = "'0.9984"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("D20").Value = "2.123.42"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("D21").Value = "'0.9982"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "'4.778"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").Value = "'5.724"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'9.077"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "'149.92"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "'141.86"
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("D27").Value = "'17.10"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "'1.918"
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("D29").Value = "'1.391"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").Value = "'0.08771"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "'3.966"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "'0.7144"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "'1.117"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "'0.01835"
$ws.Range("E37").Value = "  +8.64%  "
$ws.Range("D38").Value = "'2.675"
$ws.Range("E38").Value = "  -5.57%  "
$ws.Range("D39").Value = "'2.181"
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("D40").Value = "'0.9327"
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("D41").Value = "'5.824"
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("D42").Value = "'0.4257"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "'0.9982"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'102.24"
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("D45").Value = "'7.484"
$ws.Range("E45").Value = "  -2.19%  "
$ws.Range("D46").Value = "'0.1266"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").Value = "'32.49"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").Value = "'0.3790"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "'8.304"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("E51").Value = "  -1.43%  "

Write-Host "Updated cryptos list"